$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting the old row 4 (sCs/Nppc/Npr2/sCs) down to row 5
$ws.Rows.Item(4).Insert()

# --- Row 2 (sCs -> ECs) updated values ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.065932333333333
$ws.Range("H2").Value = 6.197797
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.33435
$ws.Range("N2").Value = 40.00305
$ws.Range("O2").Value = 0.5378090458822617
$ws.Range("P2").Value = 0.5378090458822617
$ws.Range("Q2").Value = 27.54786480898333
$ws.Range("R2").Value = 247.93078328085
$ws.Range("S2").Value = 0.5378090458822617
$ws.Range("T2").Value = 0.5378090458822617

# --- Row 3 (sCs -> FAPs) updated values ---
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.065932333333333
$ws.Range("H3").Value = 6.197797
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.523445333333334
$ws.Range("N3").Value = 16.570336
$ws.Range("O3").Value = 0.2227749282644322
$ws.Range("P3").Value = 0.2227749282644322
$ws.Range("Q3").Value = 11.41106430553244
$ws.Range("R3").Value = 102.699578749792
$ws.Range("S3").Value = 0.2227749282644322
$ws.Range("T3").Value = 0.2227749282644322

# --- Row 4 (new row: sCs -> Neutro) ---
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Nppc"
$ws.Range("C4").Value = "Npr2"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.065932333333333
$ws.Range("H4").Value = 6.197797
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05743633333333333
$ws.Range("N4").Value = 0.172309
$ws.Range("O4").Value = 0.002316556834714519
$ws.Range("P4").Value = 0.002316556834714519
$ws.Range("Q4").Value = 0.1186595781414444
$ws.Range("R4").Value = 1.067936203273
$ws.Range("S4").Value = 0.002316556834714519
$ws.Range("T4").Value = 0.002316556834714519

# --- Row 5 (previously row 4: sCs -> sCs) updated values ---
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Nppc"
$ws.Range("C5").Value = "Npr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.065932333333333
$ws.Range("H5").Value = 6.197797
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.878605666666666
$ws.Range("N5").Value = 17.635817
$ws.Range("O5").Value = 0.2370994690185916
$ws.Range("P5").Value = 0.2370994690185916
$ws.Range("Q5").Value = 12.14480152168322
$ws.Range("R5").Value = 109.303213695149
$ws.Range("S5").Value = 0.2370994690185916
$ws.Range("T5").Value = 0.2370994690185916
